# Update the cryptos list (prices / 1h volume %) with the latest scraped
# values, and fix the 47-49 coin ranking order (WEMIXToken moved up above
# NEARProtocol and Quant).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-looking values such as "1.001" or "112.90" would otherwise be
# auto-parsed by Excel as numbers (losing the original textual formatting
# such as trailing zeros), so they are entered with a leading apostrophe
# to force them to stay as text, exactly like the original cells.

$ws.Range("D2").Value = "28.024.11"
$ws.Range("D3").Value = "1.830.58"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'324.61"
$ws.Range("E5").Value = "  -2.85%  "
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "'0.4653"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").Value = "'0.3862"
$ws.Range("E8").Value = "  -1.56%  "
$ws.Range("D9").Value = "'0.07865"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").Value = "'0.9589"
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("D11").Value = "'21.85"
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("D12").Value = "1.807.01"
$ws.Range("E12").Value = "  -7.45%  "
$ws.Range("D13").Value = "'5.674"
$ws.Range("E13").Value = "  -3.03%  "
$ws.Range("D14").Value = "'6.897"
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").Value = "'0.06863"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "'87.17"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "'0.000009910"
$ws.Range("E18").Value = "  -1.63%  "
$ws.Range("D19").Value = "'16.57"
$ws.Range("E19").Value = "  -2.93%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").Value = "28.064.93"
$ws.Range("E21").Value = "  -1.72%  "
$ws.Range("D22").Value = "'5.313"
$ws.Range("E22").Value = "  -1.50%  "
$ws.Range("E23").Value = "  -2.58%  "
$ws.Range("D24").Value = "'2.092"
$ws.Range("E24").Value = "  -1.59%  "
$ws.Range("D25").Value = "2.063.76"
$ws.Range("E25").Value = "  -6.08%  "
$ws.Range("D26").Value = "'153.57"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'19.12"
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("D28").Value = "'5.710"
$ws.Range("E28").Value = "  -6.91%  "
$ws.Range("E29").Value = "  -2.98%  "
$ws.Range("D30").Value = "'117.39"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "'0.9362"
$ws.Range("E31").Value = "  -4.31%  "
$ws.Range("D32").Value = "'0.09255"
$ws.Range("E32").Value = "  -1.91%  "
$ws.Range("D33").Value = "'5.272"
$ws.Range("E33").Value = "  -1.81%  "
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("D35").Value = "'3.290"
$ws.Range("E35").Value = "  -5.94%  "
$ws.Range("D36").Value = "'0.05858"
$ws.Range("E36").Value = "  -4.40%  "
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("D38").Value = "'1.143"
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("D39").Value = "'7.779"
$ws.Range("E39").Value = "  +2.24%  "
$ws.Range("D40").Value = "'0.5579"
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("D41").Value = "'9.837"
$ws.Range("E41").Value = "  -2.64%  "
$ws.Range("E42").Value = "  -2.11%  "
$ws.Range("D43").Value = "'11.57"
$ws.Range("E43").Value = "  -2.57%  "
$ws.Range("D44").Value = "'0.07018"
$ws.Range("E44").Value = "  -1.95%  "
$ws.Range("D45").Value = "'0.5243"
$ws.Range("E45").Value = "  -2.58%  "
$ws.Range("D46").Value = "'2.124"
$ws.Range("E46").Value = "  -11.45%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'1.117"
$ws.Range("E47").Value = "  -8.72%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.824"
$ws.Range("E48").Value = "  -4.23%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'112.90"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("D50").Value = "'0.9999"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "'2.317"
$ws.Range("E51").Value = "  +0.19%  "